$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# ---------------------------------------------------------------------------
# Update the time value stored in column D (rows 2-69).
# Old value corresponded to 11:39:17 (0.48561342592592593), new value
# corresponds to 11:03:15 (0.46059027777777778) as a fraction of a day.
# ---------------------------------------------------------------------------
$newTime = 0.46059027777777778

$rng = $ws.Range("D2:D69")
$rng.Value = $newTime

# Give the updated cells an explicit black font color (rgb FF000000) instead
# of the implicit theme color, matching the new style used for this range.
$rng.Font.Color = 0

# ---------------------------------------------------------------------------
# Update the active selection on the sheet to D2:D69 with D2 as active cell.
# ---------------------------------------------------------------------------
$rng.Select() | Out-Null
